# Auto-generated script to apply scheduled market-price/profit updates
# to the Excalibur_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 14410.637
$ws.Range("J88").Value = 20817.143
$ws.Range("L88").Value = 20817.143
$ws.Range("N88").Value = -21629.143
$ws.Range("H91").Value = 14410.637
$ws.Range("J91").Value = 20817.143
$ws.Range("L91").Value = 20817.143
$ws.Range("N91").Value = -23625.143
$ws.Range("H100").Value = 6607.5386
$ws.Range("I100").Value = 2485
$ws.Range("J100").Value = 13203.6
$ws.Range("K100").Value = 2485
$ws.Range("L100").Value = 13203.6
$ws.Range("M100").Value = -1944
$ws.Range("N100").Value = -14285.6
$ws.Range("H116").Value = 190377.05
$ws.Range("J116").Value = 398145.78
$ws.Range("L116").Value = 398145.78
$ws.Range("N116").Value = -405029.78
$ws.Range("H135").Value = 624.619
$ws.Range("I135").Value = 526.4706
$ws.Range("J135").Value = 1041.75
$ws.Range("K135").Value = 4738.2354
$ws.Range("L135").Value = 9375.75
$ws.Range("M135").Value = -2203.2354
$ws.Range("N135").Value = -14445.75
$ws.Range("H137").Value = 10369.714
$ws.Range("I137").Value = 3577.3333
$ws.Range("K137").Value = 10731.9999
$ws.Range("M137").Value = -8181.999899999999
$ws.Range("H141").Value = 2956.75
$ws.Range("I141").Value = 2956.75
$ws.Range("K141").Value = 8870.25
$ws.Range("M141").Value = -3690.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21519.852
$ws.Range("I32").Value = 21831.566
$ws.Range("K32").Value = 21831.566
$ws.Range("M32").Value = -21544.566
$ws.Range("H61").Value = 742968.25
$ws.Range("I61").Value = 795849.7
$ws.Range("K61").Value = 795849.7
$ws.Range("M61").Value = -795637.7
$ws.Range("H63").Value = 8691.77
$ws.Range("I63").Value = 3998
$ws.Range("K63").Value = 3998
$ws.Range("M63").Value = -3312
$ws.Range("H66").Value = 8691.77
$ws.Range("I66").Value = 3998
$ws.Range("K66").Value = 19990
$ws.Range("M66").Value = -16558
$ws.Range("H74").Value = 3252.16
$ws.Range("I74").Value = 1657.6154
$ws.Range("J74").Value = 4979.5835
$ws.Range("K74").Value = 1657.6154
$ws.Range("L74").Value = 4979.5835
$ws.Range("M74").Value = -783.6153999999999
$ws.Range("N74").Value = -6727.5835
$ws.Range("H77").Value = 3252.16
$ws.Range("I77").Value = 1657.6154
$ws.Range("J77").Value = 4979.5835
$ws.Range("K77").Value = 8288.076999999999
$ws.Range("L77").Value = 24897.9175
$ws.Range("M77").Value = -3920.076999999999
$ws.Range("N77").Value = -33633.9175
$ws.Range("H88").Value = 2661.1667
$ws.Range("I88").Value = 2348
$ws.Range("J88").Value = 3099.6
$ws.Range("K88").Value = 2348
$ws.Range("L88").Value = 3099.6
$ws.Range("M88").Value = -1942
$ws.Range("N88").Value = -3911.6
$ws.Range("H91").Value = 2661.1667
$ws.Range("I91").Value = 2348
$ws.Range("J91").Value = 3099.6
$ws.Range("K91").Value = 2348
$ws.Range("L91").Value = 3099.6
$ws.Range("M91").Value = -944
$ws.Range("N91").Value = -5907.6
$ws.Range("H110").Value = 1483.3334
$ws.Range("I110").Value = 975
$ws.Range("J110").Value = 2500
$ws.Range("K110").Value = 975
$ws.Range("L110").Value = 2500
$ws.Range("M110").Value = 1070
$ws.Range("N110").Value = -6590
$ws.Range("H136").Value = 742968.25
$ws.Range("I136").Value = 795849.7
$ws.Range("K136").Value = 2387549.1
$ws.Range("M136").Value = -2384999.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 833.0909
$ws.Range("I94").Value = 906
$ws.Range("J94").Value = 505
$ws.Range("K94").Value = 906
$ws.Range("L94").Value = 505
$ws.Range("M94").Value = -455
$ws.Range("N94").Value = -1407
$ws.Range("H131").Value = 73332.336
$ws.Range("J131").Value = 73332.336
$ws.Range("L131").Value = 73332.336
$ws.Range("N131").Value = -83412.336
$ws.Range("H134").Value = 1137184.6
$ws.Range("I134").Value = 1083898.1
$ws.Range("K134").Value = 3251694.3
$ws.Range("M134").Value = -3249159.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 5303.273
$ws.Range("I122").Value = 1952
$ws.Range("J122").Value = 8654.546
$ws.Range("K122").Value = 5856
$ws.Range("L122").Value = 25963.638
$ws.Range("M122").Value = -3406
$ws.Range("N122").Value = -30863.638
$ws.Range("H132").Value = 202240.05
$ws.Range("I132").Value = 1773.8948
$ws.Range("K132").Value = 5321.6844
$ws.Range("M132").Value = -2791.6844

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 8460.799999999999
$ws.Range("I56").Value = 8460.799999999999
$ws.Range("K56").Value = 8460.799999999999
$ws.Range("M56").Value = -7930.799999999999
$ws.Range("H57").Value = 1336
$ws.Range("I57").Value = 1336
$ws.Range("K57").Value = 4008
$ws.Range("M57").Value = -3449
$ws.Range("H63").Value = 9843.25
$ws.Range("I63").Value = 1597.5
$ws.Range("J63").Value = 12591.833
$ws.Range("K63").Value = 4792.5
$ws.Range("L63").Value = 37775.499
$ws.Range("M63").Value = -4043.5
$ws.Range("N63").Value = -39273.499
$ws.Range("H66").Value = 9843.25
$ws.Range("I66").Value = 1597.5
$ws.Range("J66").Value = 12591.833
$ws.Range("K66").Value = 14377.5
$ws.Range("L66").Value = 113326.497
$ws.Range("M66").Value = -10633.5
$ws.Range("N66").Value = -120814.497

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4935.5386
$ws.Range("I132").Value = 3295.0625
$ws.Range("J132").Value = 7560.3
$ws.Range("K132").Value = 9885.1875
$ws.Range("L132").Value = 22680.9
$ws.Range("M132").Value = -7355.1875
$ws.Range("N132").Value = -27740.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3957.8333
$ws.Range("I7").Value = 3832.7222
$ws.Range("J7").Value = 4333.1665
$ws.Range("K7").Value = 3832.7222
$ws.Range("L7").Value = 4333.1665
$ws.Range("M7").Value = -3720.7222
$ws.Range("N7").Value = -4557.1665
$ws.Range("H23").Value = 13333.333
$ws.Range("I23").Value = 12500
$ws.Range("J23").Value = 15000
$ws.Range("K23").Value = 12500
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = -12270
$ws.Range("N23").Value = -15460
$ws.Range("H46").Value = 749.25
$ws.Range("I46").Value = 599
$ws.Range("K46").Value = 599
$ws.Range("M46").Value = -411
$ws.Range("H80").Value = 60032
$ws.Range("J80").Value = 60032
$ws.Range("L80").Value = 60032
$ws.Range("N80").Value = -62278
$ws.Range("H83").Value = 60032
$ws.Range("J83").Value = 60032
$ws.Range("L83").Value = 180096
$ws.Range("N83").Value = -191328
$ws.Range("H93").Value = 874.4
$ws.Range("I93").Value = 901.1429000000001
$ws.Range("K93").Value = 901.1429000000001
$ws.Range("M93").Value = 346.8570999999999
$ws.Range("H126").Value = 3957.8333
$ws.Range("I126").Value = 3832.7222
$ws.Range("J126").Value = 4333.1665
$ws.Range("K126").Value = 11498.1666
$ws.Range("L126").Value = 12999.4995
$ws.Range("M126").Value = -9028.1666
$ws.Range("N126").Value = -17939.4995
$ws.Range("H130").Value = 85939
$ws.Range("J130").Value = 85939
$ws.Range("L130").Value = 85939
$ws.Range("N130").Value = -95979
$ws.Range("H132").Value = 21663.076
$ws.Range("I132").Value = 26662.1
$ws.Range("K132").Value = 79986.29999999999
$ws.Range("M132").Value = -77456.29999999999
$ws.Range("H136").Value = 2582.5881
$ws.Range("I136").Value = 1858.3077
$ws.Range("J136").Value = 4936.5
$ws.Range("K136").Value = 5574.9231
$ws.Range("L136").Value = 14809.5
$ws.Range("M136").Value = -3024.9231
$ws.Range("N136").Value = -19909.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 656.6667
$ws.Range("I100").Value = 586
$ws.Range("J100").Value = 727.3333
$ws.Range("K100").Value = 1172
$ws.Range("L100").Value = 1454.6666
$ws.Range("M100").Value = -631
$ws.Range("N100").Value = -2536.6666
$ws.Range("H122").Value = 1789.1724
$ws.Range("I122").Value = 1467.8182
$ws.Range("J122").Value = 2799.1428
$ws.Range("K122").Value = 4403.4546
$ws.Range("L122").Value = 8397.428400000001
$ws.Range("M122").Value = -1953.4546
$ws.Range("N122").Value = -13297.4284
$ws.Range("H136").Value = 6454.9023
$ws.Range("I136").Value = 6653
$ws.Range("J136").Value = 5492.7144
$ws.Range("K136").Value = 19959
$ws.Range("L136").Value = 16478.1432
$ws.Range("M136").Value = -17409
